# Unit 2 - LC 101 - Class 1: slide 2 ("Expectations" shape, id 75)
# Split the "Expectations:" paragraph so a new, blank, bullet-less
# paragraph follows it (simulates pressing Enter right after the colon),
# then re-point the per-paragraph click animations that targeted the
# paragraphs pushed down by the new line.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)              # Shape 75
$tr = $sh.TextFrame.TextRange

# Paragraph 5 (1-based) is "Expectations:" - break a new paragraph
# right after it.
$expectationsPara = $tr.Paragraphs(5, 1)
$null = $expectationsPara.InsertAfter("`r")

# The newly created paragraph is now #6; it inherited the bullet from
# "Expectations:" - this blank continuation line shouldn't show one.
$newPara = $tr.Paragraphs(6, 1)
$newPara.ParagraphFormat.Bullet.Visible = 0

# The click-animation effects that targeted paragraphs after the split
# point (old 1-based indices 6, 7, 9, 10) need to move down one slot.
$ms = $s.TimeLine.MainSequence
for ($i = 1; $i -le $ms.Count; $i++) {
    $eff = $ms.Item($i)
    if ($eff.Shape.Id -eq 75 -and $eff.Paragraph -ge 6) {
        $eff.Paragraph = $eff.Paragraph + 1
    }
}
